$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
  2 = @("ECs", "Lta", "Ltbr", "ECs", 1, 0.3333333333333333, 0.103879, 0.311637, 0.1899090787212519, 0.1899090787212519, 3, 1, 12.87437866666667, 38.623136, 0.272778495601419, 0.272778495601419, 1.337377581514667, 12.036398233632, 0.05180311279463455, 0.05180311279463455)
  3 = @("ECs", "Lta", "Ltbr", "FAPs", 1, 0.3333333333333333, 0.103879, 0.311637, 0.1899090787212519, 0.1899090787212519, 3, 1, 17.704387, 53.11316100000001, 0.3751152716914535, 0.3751152716914535, 1.839114017173, 16.552026154557, 0.07123779566119605, 0.07123779566119605)
  4 = @("ECs", "Lta", "Ltbr", "MuSCs", 1, 0.3333333333333333, 0.103879, 0.311637, 0.1899090787212519, 0.1899090787212519, 3, 1, 8.077278666666666, 24.231836, 0.1711389714636405, 0.1711389714636405, 0.8390596306146666, 7.551536675532001, 0.03250084440396259, 0.03250084440396259)
  5 = @("ECs", "Lta", "Ltbr", "Resolving-Mac", 1, 0.3333333333333333, 0.103879, 0.311637, 0.1899090787212519, 0.1899090787212519, 3, 1, 8.541146333333332, 25.623439, 0.180967261243487, 0.180967261243487, 0.8872457399603332, 7.985211659642999, 0.03436732586145873, 0.03436732586145873)
  6 = @("FAPs", "Lta", "Ltbr", "ECs", 2, 0.6666666666666666, 0.257314, 0.771942, 0.4704152396738534, 0.4704152396738534, 3, 1, 12.87437866666667, 38.623136, 0.272778495601419, 0.272778495601419, 3.312757872234667, 29.814820850112, 0.1283191613862147, 0.1283191613862147)
  7 = @("FAPs", "Lta", "Ltbr", "FAPs", 2, 0.6666666666666666, 0.257314, 0.771942, 0.4704152396738534, 0.4704152396738534, 3, 1, 17.704387, 53.11316100000001, 0.3751152716914535, 0.3751152716914535, 4.555586636518, 41.00027972866201, 0.1764599404380577, 0.1764599404380577)
  8 = @("FAPs", "Lta", "Ltbr", "MuSCs", 2, 0.6666666666666666, 0.257314, 0.771942, 0.4704152396738534, 0.4704152396738534, 3, 1, 8.077278666666666, 24.231836, 0.1711389714636405, 0.1711389714636405, 2.078396882834666, 18.705571945512, 0.08050638027860521, 0.0805063802786052)
  9 = @("FAPs", "Lta", "Ltbr", "Resolving-Mac", 2, 0.6666666666666666, 0.257314, 0.771942, 0.4704152396738534, 0.4704152396738534, 3, 1, 8.541146333333332, 25.623439, 0.180967261243487, 0.180967261243487, 2.197756527615333, 19.779808748538, 0.08512975757097578, 0.08512975757097577)
  10 = @("Resolving-Mac", "Lta", "Ltbr", "ECs", 2, 0.6666666666666666, 0.1858003333333333, 0.557401, 0.3396756816048946, 0.3396756816048946, 3, 1, 12.87437866666667, 38.623136, 0.272778495601419, 0.272778495601419, 2.392063847726222, 21.528574629536, 0.09265622142056974, 0.09265622142056974)
  11 = @("Resolving-Mac", "Lta", "Ltbr", "FAPs", 2, 0.6666666666666666, 0.1858003333333333, 0.557401, 0.3396756816048946, 0.3396756816048946, 3, 1, 17.704387, 53.11316100000001, 0.3751152716914535, 0.3751152716914535, 3.289481006062334, 29.605329054561, 0.1274175355921997, 0.1274175355921997)
  12 = @("Resolving-Mac", "Lta", "Ltbr", "MuSCs", 2, 0.6666666666666666, 0.1858003333333333, 0.557401, 0.3396756816048946, 0.3396756816048946, 3, 1, 8.077278666666666, 24.231836, 0.1711389714636405, 0.1711389714636405, 1.500761068692889, 13.506849618236, 0.0581317467810727, 0.0581317467810727)
  13 = @("Resolving-Mac", "Lta", "Ltbr", "Resolving-Mac", 2, 0.6666666666666666, 0.1858003333333333, 0.557401, 0.3396756816048946, 0.3396756816048946, 3, 1, 8.541146333333332, 25.623439, 0.180967261243487, 0.180967261243487, 1.586947835782111, 14.282530522039, 0.06147017781105248, 0.06147017781105247)
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item([int]$r, $i + 1).Value = $vals[$i]
  }
}
